# Update scripts with new TPM values.
# The "Sending cluster" labels are renamed:
#   MuSCs          -> Inflammatory-Mac  (rows 2-4)
#   Resolving-Mac  -> MuSCs             (rows 5-7)
# and the derived NATMI metrics (columns E..T) are recomputed for the
# new TPM-based clustering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sending-cluster labels (column A) ---
$ws.Range("A2:A4").Value = "Inflammatory-Mac"
$ws.Range("A5:A7").Value = "MuSCs"

# --- Row 2 (Inflammatory-Mac -> ECs) ---
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05008433333333333
$ws.Range("H2").Value = 0.150253
$ws.Range("I2").Value = 0.3054767171413236
$ws.Range("J2").Value = 0.3054767171413236
$ws.Range("M2").Value = 0.106124
$ws.Range("N2").Value = 0.318372
$ws.Range("O2").Value = 0.08094716512538251
$ws.Range("P2").Value = 0.08094716512538253
$ws.Range("Q2").Value = 0.005315149790666667
$ws.Range("R2").Value = 0.047836348116
$ws.Range("S2").Value = 0.02472747426439849
$ws.Range("T2").Value = 0.0247274742643985

# --- Row 3 (Inflammatory-Mac -> FAPs) ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05008433333333333
$ws.Range("H3").Value = 0.150253
$ws.Range("I3").Value = 0.3054767171413236
$ws.Range("J3").Value = 0.3054767171413236
$ws.Range("O3").Value = 0.8331551016962769
$ws.Range("P3").Value = 0.833155101696277
$ws.Range("Q3").Value = 0.05470659976188889
$ws.Range("R3").Value = 0.492359397857
$ws.Range("S3").Value = 0.2545094853357243
$ws.Range("T3").Value = 0.2545094853357243

# --- Row 4 (Inflammatory-Mac -> MuSCs) ---
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05008433333333333
$ws.Range("H4").Value = 0.150253
$ws.Range("I4").Value = 0.3054767171413236
$ws.Range("J4").Value = 0.3054767171413236
$ws.Range("M4").Value = 0.1126143333333333
$ws.Range("N4").Value = 0.337843
$ws.Range("O4").Value = 0.08589773317834044
$ws.Range("P4").Value = 0.08589773317834046
$ws.Range("Q4").Value = 0.005640213808777778
$ws.Range("R4").Value = 0.050761924279
$ws.Range("S4").Value = 0.02623975754120079
$ws.Range("T4").Value = 0.0262397575412008

# --- Row 5 (MuSCs -> ECs) ---
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1138703333333333
$ws.Range("H5").Value = 0.341611
$ws.Range("I5").Value = 0.6945232828586764
$ws.Range("J5").Value = 0.6945232828586764
$ws.Range("M5").Value = 0.106124
$ws.Range("N5").Value = 0.318372
$ws.Range("O5").Value = 0.08094716512538251
$ws.Range("P5").Value = 0.08094716512538253
$ws.Range("Q5").Value = 0.01208437525466667
$ws.Range("R5").Value = 0.108759377292
$ws.Range("S5").Value = 0.05621969086098402
$ws.Range("T5").Value = 0.05621969086098403

# --- Row 6 (MuSCs -> FAPs) ---
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1138703333333333
$ws.Range("H6").Value = 0.341611
$ws.Range("I6").Value = 0.6945232828586764
$ws.Range("J6").Value = 0.6945232828586764
$ws.Range("O6").Value = 0.8331551016962769
$ws.Range("P6").Value = 0.833155101696277
$ws.Range("Q6").Value = 0.1243793884398889
$ws.Range("R6").Value = 1.119414495959
$ws.Range("S6").Value = 0.5786456163605527
$ws.Range("T6").Value = 0.5786456163605527

# --- Row 7 (MuSCs -> MuSCs) ---
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1138703333333333
$ws.Range("H7").Value = 0.341611
$ws.Range("I7").Value = 0.6945232828586764
$ws.Range("J7").Value = 0.6945232828586764
$ws.Range("M7").Value = 0.1126143333333333
$ws.Range("N7").Value = 0.337843
$ws.Range("O7").Value = 0.08589773317834044
$ws.Range("P7").Value = 0.08589773317834046
$ws.Range("Q7").Value = 0.01282343167477778
$ws.Range("R7").Value = 0.115410885073
$ws.Range("S7").Value = 0.05965797563713965
$ws.Range("T7").Value = 0.05965797563713966
